$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from an existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF), rows 2-40
$iValues = @(8,8,1,5,9,7,5,7,7,8,7,6,7,8,5,9,7,5,6,8,7,9,7,8,6,1,6,10,6,9,3,6,6,6,9,9,3,3,9)
$jValues = @(8,8,1,6,9,8,6,7,8,8,7,6,8,8,6,9,8,6,7,8,7,9,7,8,6,2,8,10,6,9,5,7,7,6,9,9,3,3,9)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}

$wb.Save()
